$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 741.5
$ws.Range("I32").Value = 700
$ws.Range("J32").Value = 766.4
$ws.Range("K32").Value = 700
$ws.Range("L32").Value = 766.4
$ws.Range("M32").Value = -374
$ws.Range("N32").Value = -1418.4

$ws.Range("H69").Value = 4730
$ws.Range("I69").Value = 5000
$ws.Range("J69").Value = 4460
$ws.Range("K69").Value = 15000
$ws.Range("L69").Value = 13380
$ws.Range("M69").Value = -14126
$ws.Range("N69").Value = -15128

$ws.Range("H72").Value = 4730
$ws.Range("I72").Value = 5000
$ws.Range("J72").Value = 4460
$ws.Range("K72").Value = 45000
$ws.Range("L72").Value = 40140
$ws.Range("M72").Value = -40632
$ws.Range("N72").Value = -48876

$ws.Range("H129").Value = 2675.6606
$ws.Range("J129").Value = 939.25
$ws.Range("L129").Value = 2817.75
$ws.Range("N129").Value = -12817.75

$ws.Range("H138").Value = 4023.1067
$ws.Range("J138").Value = 4024.7231
$ws.Range("L138").Value = 12074.1693
$ws.Range("N138").Value = -22354.1693

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 10981.3
$ws.Range("I28").Value = 10981.3
$ws.Range("K28").Value = 10981.3
$ws.Range("M28").Value = -10789.3

$ws.Range("H44").Value = 12771.125
$ws.Range("J44").Value = 12738.429
$ws.Range("L44").Value = 12738.429
$ws.Range("N44").Value = -13714.429

$ws.Range("H99").Value = 10981.3
$ws.Range("I99").Value = 10981.3
$ws.Range("K99").Value = 10981.3
$ws.Range("M99").Value = -7986.299999999999

$ws.Range("H122").Value = 2717.3333
$ws.Range("I122").Value = 2373.3333
$ws.Range("J122").Value = 4093.3333
$ws.Range("K122").Value = 7119.999899999999
$ws.Range("L122").Value = 12279.9999
$ws.Range("M122").Value = -4669.999899999999
$ws.Range("N122").Value = -17179.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 299.5
$ws.Range("I8").Value = 299.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 299.5
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -159.5
$ws.Range("N8").ClearContents()

$ws.Range("H82").Value = 17918.857
$ws.Range("I82").Value = 3164
$ws.Range("K82").Value = 3164
$ws.Range("M82").Value = -2781

$ws.Range("H85").Value = 17918.857
$ws.Range("I85").Value = 3164
$ws.Range("K85").Value = 3164
$ws.Range("M85").Value = -1838

$ws.Range("H107").Value = 100048376
$ws.Range("I107").Value = 166743520
$ws.Range("J107").Value = 5665.25
$ws.Range("K107").Value = 166743520
$ws.Range("L107").Value = 5665.25
$ws.Range("M107").Value = -166741600
$ws.Range("N107").Value = -9505.25

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 958.94446
$ws.Range("I16").Value = 704.36365
$ws.Range("J16").Value = 1359
$ws.Range("K16").Value = 704.36365
$ws.Range("L16").Value = 1359
$ws.Range("M16").Value = -417.36365
$ws.Range("N16").Value = -1933

$ws.Range("H31").Value = 142262.6
$ws.Range("I31").Value = 1576.5
$ws.Range("J31").Value = 705007
$ws.Range("K31").Value = 1576.5
$ws.Range("L31").Value = 705007
$ws.Range("M31").Value = -1281.5
$ws.Range("N31").Value = -705597

$ws.Range("H34").Value = 142262.6
$ws.Range("I34").Value = 1576.5
$ws.Range("J34").Value = 705007
$ws.Range("K34").Value = 1576.5
$ws.Range("L34").Value = 705007
$ws.Range("M34").Value = -1374.5
$ws.Range("N34").Value = -705411

$ws.Range("H62").Value = 2662.6365
$ws.Range("I62").Value = 2495
$ws.Range("J62").Value = 2699.889
$ws.Range("K62").Value = 2495
$ws.Range("L62").Value = 2699.889
$ws.Range("M62").Value = -1871
$ws.Range("N62").Value = -3947.889

$ws.Range("H65").Value = 2662.6365
$ws.Range("I65").Value = 2495
$ws.Range("J65").Value = 2699.889
$ws.Range("K65").Value = 12475
$ws.Range("L65").Value = 13499.445
$ws.Range("M65").Value = -9355
$ws.Range("N65").Value = -19739.445

$ws.Range("H99").Value = 12539.4
$ws.Range("I99").Value = 3193.3333
$ws.Range("J99").Value = 16544.857
$ws.Range("K99").Value = 3193.3333
$ws.Range("L99").Value = 16544.857
$ws.Range("M99").Value = -1695.3333
$ws.Range("N99").Value = -19540.857

$ws.Range("H113").Value = 958.94446
$ws.Range("I113").Value = 704.36365
$ws.Range("J113").Value = 1359
$ws.Range("K113").Value = 704.36365
$ws.Range("L113").Value = 1359
$ws.Range("M113").Value = 1465.63635
$ws.Range("N113").Value = -5699

$ws.Range("H115").Value = 39999
$ws.Range("J115").Value = 39999
$ws.Range("L115").Value = 39999
$ws.Range("N115").Value = -42349

$ws.Range("H123").Value = 45000
$ws.Range("J123").Value = 45000
$ws.Range("L123").Value = 45000
$ws.Range("N123").Value = -54800

$ws.Range("H126").Value = 12539.4
$ws.Range("I126").Value = 3193.3333
$ws.Range("J126").Value = 16544.857
$ws.Range("K126").Value = 9579.999899999999
$ws.Range("L126").Value = 49634.571
$ws.Range("M126").Value = -7109.999899999999
$ws.Range("N126").Value = -54574.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 68.21429000000001
$ws.Range("I12").Value = 30.5
$ws.Range("J12").Value = 74.5
$ws.Range("K12").Value = 91.5
$ws.Range("L12").Value = 223.5
$ws.Range("M12").Value = 81.5
$ws.Range("N12").Value = -569.5

$ws.Range("H47").Value = 353
$ws.Range("I47").Value = 346.66666
$ws.Range("J47").Value = 372
$ws.Range("K47").Value = 1039.99998
$ws.Range("L47").Value = 1116
$ws.Range("M47").Value = -608.9999800000001
$ws.Range("N47").Value = -1978

$ws.Range("H131").Value = 714901.6
$ws.Range("J131").Value = 746407.5600000001
$ws.Range("L131").Value = 2239222.68
$ws.Range("N131").Value = -2249302.68

$ws.Range("H132").Value = 2431.6667
$ws.Range("I132").Value = 1822
$ws.Range("J132").Value = 2666.1538
$ws.Range("K132").Value = 16398
$ws.Range("L132").Value = 23995.3842
$ws.Range("M132").Value = -13868
$ws.Range("N132").Value = -29055.3842

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1443571.4
$ws.Range("I107").Value = 598
$ws.Range("J107").Value = 5051005
$ws.Range("K107").Value = 598
$ws.Range("L107").Value = 5051005
$ws.Range("M107").Value = 1322
$ws.Range("N107").Value = -5054845

$ws.Range("H122").Value = 2683
$ws.Range("I122").Value = 2428.5
$ws.Range("J122").Value = 3192
$ws.Range("K122").Value = 7285.5
$ws.Range("L122").Value = 9576
$ws.Range("M122").Value = -4835.5
$ws.Range("N122").Value = -14476

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 58014.277
$ws.Range("I40").Value = 168849.5
$ws.Range("J40").Value = 2596.6667
$ws.Range("K40").Value = 168849.5
$ws.Range("L40").Value = 2596.6667
$ws.Range("M40").Value = -168713.5
$ws.Range("N40").Value = -2868.6667

$ws.Range("H46").Value = 5908.8887
$ws.Range("J46").Value = 4928.5713
$ws.Range("L46").Value = 4928.5713
$ws.Range("N46").Value = -5304.5713

$ws.Range("H93").Value = 4661
$ws.Range("I93").Value = 4661
$ws.Range("K93").Value = 4661
$ws.Range("M93").Value = -3413

$ws.Range("H122").Value = 5355
$ws.Range("I122").Value = 4734.3335
$ws.Range("K122").Value = 14203.0005
$ws.Range("M122").Value = -11753.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 35200
$ws.Range("J75").Value = 35200
$ws.Range("L75").Value = 35200
$ws.Range("N75").Value = -37072

$ws.Range("H78").Value = 35200
$ws.Range("J78").Value = 35200
$ws.Range("L78").Value = 105600
$ws.Range("N78").Value = -114960

$ws.Range("H81").Value = 250649.88
$ws.Range("I81").Value = 200720.2
$ws.Range("J81").Value = 333866
$ws.Range("K81").Value = 401440.4
$ws.Range("L81").Value = 667732
$ws.Range("M81").Value = -400379.4
$ws.Range("N81").Value = -669854

$ws.Range("H84").Value = 250649.88
$ws.Range("I84").Value = 200720.2
$ws.Range("J84").Value = 333866
$ws.Range("K84").Value = 2007202
$ws.Range("L84").Value = 3338660
$ws.Range("M84").Value = -2001898
$ws.Range("N84").Value = -3349268

$ws.Range("H126").Value = 1015.6429
$ws.Range("I126").Value = 1058.2222
$ws.Range("J126").Value = 939
$ws.Range("K126").Value = 3174.6666
$ws.Range("L126").Value = 2817
$ws.Range("M126").Value = -704.6665999999996
$ws.Range("N126").Value = -7757
